$wb = $excel.ActiveWorkbook

# ---- Rushing sheet ----
$rushing = $wb.Worksheets.Item("Rushing")

# B.Mayfield (row 2): 3DATT 5 -> 6
$rushing.Range("E2").Value = 6

# N.Chubb (row 4): 1DATT 96 -> 109, 2DATT 60 -> 64, RZATT 29 -> 31
$rushing.Range("C4").Value = 109
$rushing.Range("D4").Value = 64
$rushing.Range("F4").Value = 31

# D.Johnson (row 6): 1DATT 5 -> 7, 2DATT 4 -> 5, 3DATT 3 -> 4, RZATT 0 -> 1
$rushing.Range("C6").Value = 7
$rushing.Range("D6").Value = 5
$rushing.Range("E6").Value = 4
$rushing.Range("F6").Value = 1

# A.Schwartz (row 12): 1DATT 1 -> 3
$rushing.Range("C12").Value = 3

# ---- Receiving sheet ----
$receiving = $wb.Worksheets.Item("Receiving")

# N.Chubb (row 2): Short Target 12 -> 15, Short Comp 9 -> 11, Deep Target 1 -> 2, Deep Comp 1 -> 2
$receiving.Range("C2").Value = 15
$receiving.Range("D2").Value = 11
$receiving.Range("E2").Value = 2
$receiving.Range("F2").Value = 2

# D.Johnson (row 4): Short Target 5 -> 7, Short Comp 5 -> 6
$receiving.Range("C4").Value = 7
$receiving.Range("D4").Value = 6

# A.Janovich (row 6): Short Target 1 -> 3, Short Comp 1 -> 3
$receiving.Range("C6").Value = 3
$receiving.Range("D6").Value = 3

# J.Landry (row 8): Short Target 58 -> 63, Short Comp 47 -> 50, Deep Target 22 -> 24,
#                   Deep Comp 11 -> 12, RZ Target 7 -> 8, RZ Comp 3 -> 4
$receiving.Range("C8").Value = 63
$receiving.Range("D8").Value = 50
$receiving.Range("E8").Value = 24
$receiving.Range("F8").Value = 12
$receiving.Range("G8").Value = 8
$receiving.Range("H8").Value = 4

# D.Peoples-Jones (row 9): Short Target 20 -> 23, Short Comp 15 -> 16, Deep Target 15 -> 18
$receiving.Range("C9").Value = 23
$receiving.Range("D9").Value = 16
$receiving.Range("E9").Value = 18

# R.Higgins (row 10): Short Target 24 -> 30, Short Comp 14 -> 19, RZ Target 3 -> 5, RZ Comp 2 -> 3
$receiving.Range("C10").Value = 30
$receiving.Range("D10").Value = 19
$receiving.Range("G10").Value = 5
$receiving.Range("H10").Value = 3

# A.Schwartz (row 11): Short Target 9 -> 10, Short Comp 7 -> 8, RZ Target 0 -> 1, RZ Comp 0 -> 1
$receiving.Range("C11").Value = 10
$receiving.Range("D11").Value = 8
$receiving.Range("G11").Value = 1
$receiving.Range("H11").Value = 1

# A.Hooper (row 13): Short Target 48 -> 52, Short Comp 30 -> 33, Deep Target 5 -> 6,
#                    RZ Target 10 -> 11, RZ Comp 6 -> 7
$receiving.Range("C13").Value = 52
$receiving.Range("D13").Value = 33
$receiving.Range("E13").Value = 6
$receiving.Range("G13").Value = 11
$receiving.Range("H13").Value = 7

# D.Njoku (row 14): Short Target 27 -> 29
$receiving.Range("C14").Value = 29

# H.Bryant (row 15): Short Target 15 -> 16, Short Comp 12 -> 13, RZ Target 1 -> 2, RZ Comp 1 -> 2
$receiving.Range("C15").Value = 16
$receiving.Range("D15").Value = 13
$receiving.Range("G15").Value = 2
$receiving.Range("H15").Value = 2
